$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.354.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -2.65%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.184.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -4.09%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.04%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''585.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -2.50%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''135.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -6.35%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.03%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''3.182.05'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -4.16%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -4.20%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.141'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -5.81%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''5.24'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -5.93%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.450'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -5.38%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.0000234'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -6.62%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''33.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -5.18%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.709.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -4.07%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = '''  -1.50%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''3.185.64'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -4.12%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''62.410.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -2.66%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''6.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -4.97%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''456.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -5.60%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''13.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -3.12%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.703'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -4.99%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''7.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -5.34%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''13.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -2.63%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''82.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -3.02%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.24%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = 'FirstDigitalUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -0.03%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '''2.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -4.07%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''6.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -5.11%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''7.81'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -6.09%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -7.67%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -8.32%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  -3.51%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  -7.58%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -6.45%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''5.79'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -3.75%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -4.09%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -10.15%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.0387'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -4.32%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''410.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -5.83%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''2.941.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -3.88%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  +0.65%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''8.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -5.40%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''2.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -6.63%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -3.94%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  -7.31%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''35.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Value = '''25.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -4.35%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''122.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -0.68%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -4.58%  '
$ws.Range("E51").Style = "Normal"
